{"js": "// Fix exercise 2 boundary-value test data: update the \"Cadastrar produto no\n// valor de R$ X\" entries in the Exerc\u00edcio 2 test-case table so the\n// boundary values actually sit just inside/outside the valid range\n// (R$19,00 .. R$99,00) instead of on the round numbers.\n\nconst replacements = [\n  { from: \"Cadastrar produto no valor de R$ 18,00\", to: \"Cadastrar produto no valor de R$ 18,99\" },\n  { from: \"Cadastrar produto no valor de R$ 20,00\", to: \"Cadastrar produto no valor de R$ 19,01\" },\n  { from: \"Cadastrar produto no valor de R$ 98,00\", to: \"Cadastrar produto no valor de R$ 98,99\" },\n  { from: \"Cadastrar produto no valor de R$ 100,00\", to: \"Cadastrar produto no valor de R$ 99,01\" },\n];\n\n// These edits all live in the second exercise's test-case table (the one\n// right after the \"Exerc\u00edcio 2\" heading). Scope the search to that table so\n// we don't touch the similarly-worded rows that belong to other exercises\n// (e.g. the \"R$ 100,00\" row in exercise 1's table must stay untouched).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst exercise2Table = tables.items[3];\n\nfor (const { from, to } of replacements) {\n  const results = exercise2Table.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Fix exercise 2 boundary-value test data: update the \"Cadastrar produto no\n# valor de R$ X\" entries in the Exerc\u00edcio 2 test-case table so the boundary\n# values actually sit just inside/outside the valid range (R$19,00 .. R$99,00)\n# instead of landing on the round numbers.\n\n$d = $word.ActiveDocument\n\n# The table right after the \"Exerc\u00edcio 2\" heading is the 4th table in the\n# document (1-based COM indexing). Scope every Find/Replace to that table's\n# Range so the similarly-worded rows belonging to other exercises (e.g. the\n# \"R$ 100,00\" row that legitimately belongs to exercise 1's table) are left\n# untouched.\n$exercise2Table = $d.Tables.Item(4)\n\nfunction Replace-InRange($table, [string]$findText, [string]$replaceText) {\n    $rng = $table.Range\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    # wdReplaceOne (1) so only the single match inside this table's range is\n    # touched -- wdReplaceAll (2) ends up replacing every match in the whole\n    # document, which would also clobber exercise 1's \"R$ 100,00\" row.\n    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 1)\n}\n\nReplace-InRange $exercise2Table \"Cadastrar produto no valor de R$ 18,00\" \"Cadastrar produto no valor de R$ 18,99\"\nReplace-InRange $exercise2Table \"Cadastrar produto no valor de R$ 20,00\" \"Cadastrar produto no valor de R$ 19,01\"\nReplace-InRange $exercise2Table \"Cadastrar produto no valor de R$ 98,00\" \"Cadastrar produto no valor de R$ 98,99\"\nReplace-InRange $exercise2Table \"Cadastrar produto no valor de R$ 100,00\" \"Cadastrar produto no valor de R$ 99,01\"\n"}
